# "Generate Report for Handoff"
#
# This updates the localization-status report to reflect that the
# handoff package is now ready (instead of showing the previous
# handback status), refreshes the associated timestamps, and shrinks
# the "Status"/"zh-cn"/"de-de" columns now that the status text is
# shorter than before.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)   # "Overview" sheet
$zhcn     = $wb.Worksheets.Item(2)   # "zh-cn" sheet
$dede     = $wb.Worksheets.Item(3)   # "de-de" sheet

# --- Update status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
# Overview sheet shows this status for both locales (columns E and F)
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
# Per-locale sheets show it in the "Status" column (column C)
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Refresh timestamps ---
# Overview "Latest HO Xliff Generate Date" (column G), and de-de's
# "Latest Handoff Datetime" (column H), which happen to share the same
# underlying timestamp value as each other.
$overview.Range("G2").Value = "2016-09-05 07:09:29"
$dede.Range("H2").Value = "2016-09-05 07:09:29"
# zh-cn "Latest Handoff Datetime" (column H) has its own distinct timestamp
$zhcn.Range("H2").Value = "2016-09-05 07:09:24"

# --- Shrink the now-narrower status columns to fit the shorter text ---
$overview.Range("E1").ColumnWidth = 16.3
$overview.Range("F1").ColumnWidth = 16.3
$zhcn.Range("C1").ColumnWidth = 16.3
$dede.Range("C1").ColumnWidth = 16.3
